# Saldo_guide.xlsx update
#  - Rename sheet (new export timestamp in the name)
#  - Refresh the "Dt. Referencia" column (G) for every data row to the new date
#  - Update a handful of balances that changed between the two exports
#  - Leave the active selection on I16, matching the new saved view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new export timestamp
$ws.Name = "IClientBalance-20241007-093313-"

# The data rows run from row 2 to row 274; refresh the reference date (column G)
# for all of them to the new value (serial date 45572 = 2024-10-07).
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45572
}

# A few account balances were refreshed between the two exports.
# Row 108 (B108 = 378873): projected value zeroed out, predicted balance absorbs it.
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 70476.36

# Row 112 (B112 = 379727): predicted balance / total updated.
$ws.Range("E112").Value = 62.01
$ws.Range("H112").Value = 62.01

# Row 161 (B161 = 445690): predicted balance / total updated.
$ws.Range("E161").Value = 298.58999999999997
$ws.Range("H161").Value = 298.58999999999997

# Row 255 (B255 = 806458): predicted balance / total updated.
$ws.Range("E255").Value = 37967.75
$ws.Range("H255").Value = 37967.75

# Restore the saved selection (I16) recorded in the sheet view.
$ws.Range("I16").Select()
